$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff" — a fresh handoff was produced for the
# 8a403151-f11e-4e7e-8714-4d1a1e62c038 file, so its "Latest Handoff
# Datetime" (column D, row 4) is refreshed on both the zh-cn and de-de
# status sheets.

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-19 06:10:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-19 06:10:13"
